# Add the PF/1.0.3 release row to the meta-sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.3"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# Keep the new row on the sheet's default (unstyled) formatting, matching
# the rest of the workbook's plain "Normal" style.
$ws.Range("A3:D3").Style = "Normal"
